# Updated cryptos list on Wed Jul  5 15:49:31 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell used to push "Price" column updates through Excel as literal
# text (PasteSpecial -Values only) so a value like "1.000" or "237.81" is not
# auto-coerced into a Double the way a direct .Value assignment would.
$helper = $ws.Cells.Item(1, 10)
function Set-PriceText($row, $val) {
    $escaped = $val.Replace('"', '""')
    $helper.Formula = '="' + $escaped + '"'
    $helper.Copy()
    $ws.Cells.Item($row, 4).PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

# Row 2: Bitcoin
Set-PriceText 2 '30.378.91'
$ws.Cells.Item(2, 5).Value = '  -1.92%  '

# Row 3: Ethereum
Set-PriceText 3 '1.904.42'
$ws.Cells.Item(3, 5).Value = '  -2.60%  '

# Row 4: TetherUSD
Set-PriceText 4 '1.000'
$ws.Cells.Item(4, 5).Value = '  -0.12%  '

# Row 5: BNB
Set-PriceText 5 '237.81'
$ws.Cells.Item(5, 5).Value = '  -2.52%  '

# Row 6: USDC
Set-PriceText 6 '1.001'
$ws.Cells.Item(6, 5).Value = '  -0.03%  '

# Row 7: XRP
Set-PriceText 7 '0.4719'
$ws.Cells.Item(7, 5).Value = '  -2.42%  '

# Row 8: Cardano
Set-PriceText 8 '0.2821'
$ws.Cells.Item(8, 5).Value = '  -4.01%  '

# Row 9: Dogecoin
Set-PriceText 9 '0.06610'
$ws.Cells.Item(9, 5).Value = '  -6.50%  '

# Row 10: Solana
Set-PriceText 10 '18.60'
$ws.Cells.Item(10, 5).Value = '  -5.52%  '

# Row 11: Litecoin
Set-PriceText 11 '99.43'
$ws.Cells.Item(11, 5).Value = '  -7.18%  '

# Row 12: TRON
Set-PriceText 12 '0.07707'
$ws.Cells.Item(12, 5).Value = '  -1.13%  '

# Row 13: WrappedEther
Set-PriceText 13 '1.906.16'
$ws.Cells.Item(13, 5).Value = '  -2.64%  '

# Row 14: Polkadot
Set-PriceText 14 '5.147'
$ws.Cells.Item(14, 5).Value = '  -5.38%  '

# Row 15: Polygon
Set-PriceText 15 '0.6635'
$ws.Cells.Item(15, 5).Value = '  -5.16%  '

# Row 16: WrappedBTC
Set-PriceText 16 '30.411.08'
$ws.Cells.Item(16, 5).Value = '  -1.84%  '

# Row 17: BitcoinCash
Set-PriceText 17 '252.69'
$ws.Cells.Item(17, 5).Value = '  -9.66%  '

# Row 18: Dai
Set-PriceText 18 '1.001'
$ws.Cells.Item(18, 5).Value = '  +0.06%  '

# Row 19: ShibaInu
Set-PriceText 19 '0.000007408'
$ws.Cells.Item(19, 5).Value = '  -4.96%  '

# Row 20: Avalanche
Set-PriceText 20 '12.59'
$ws.Cells.Item(20, 5).Value = '  -5.15%  '

# Row 21: Uniswap
Set-PriceText 21 '5.340'
$ws.Cells.Item(21, 5).Value = '  -3.61%  '

# Row 22: BinanceUSD
Set-PriceText 22 '1.000'
$ws.Cells.Item(22, 5).Value = '  -0.11%  '

# Row 23: BitDAO
$ws.Cells.Item(23, 2).Value = 'BitDAO'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
Set-PriceText 23 '0.4565'
$ws.Cells.Item(23, 5).Value = '  -7.15%  '

# Row 24: Chainlink
$ws.Cells.Item(24, 2).Value = 'Chainlink'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-PriceText 24 '6.258'
$ws.Cells.Item(24, 5).Value = '  -3.55%  '

# Row 25: Cosmos
$ws.Cells.Item(25, 2).Value = 'Cosmos'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-PriceText 25 '9.306'
$ws.Cells.Item(25, 5).Value = '  -5.11%  '

# Row 26: Monero
$ws.Cells.Item(26, 2).Value = 'Monero'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-PriceText 26 '164.88'
$ws.Cells.Item(26, 5).Value = '  -2.25%  '

# Row 27: EthereumClassic
$ws.Cells.Item(27, 2).Value = 'EthereumClassic'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-PriceText 27 '18.76'
$ws.Cells.Item(27, 5).Value = '  -5.07%  '

# Row 28: LidoDAOToken
$ws.Cells.Item(28, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-PriceText 28 '2.037'
$ws.Cells.Item(28, 5).Value = '  -6.22%  '

# Row 29: Stellar
$ws.Cells.Item(29, 2).Value = 'Stellar'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-PriceText 29 '0.1007'
$ws.Cells.Item(29, 5).Value = '  -3.79%  '

# Row 30: Toncoin
$ws.Cells.Item(30, 2).Value = 'Toncoin'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-PriceText 30 '1.378'
$ws.Cells.Item(30, 5).Value = '  -0.51%  '

# Row 31: Filecoin
$ws.Cells.Item(31, 2).Value = 'Filecoin'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-PriceText 31 '4.612'
$ws.Cells.Item(31, 5).Value = '  +0.44%  '

# Row 32: PancakeSwap
$ws.Cells.Item(32, 2).Value = 'PancakeSwap'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-PriceText 32 '1.506'
$ws.Cells.Item(32, 5).Value = '  -4.05%  '

# Row 33: InternetComputer(DFINITY)
$ws.Cells.Item(33, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-PriceText 33 '4.208'
$ws.Cells.Item(33, 5).Value = '  -5.09%  '

# Row 34: Hedera
$ws.Cells.Item(34, 2).Value = 'Hedera'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-PriceText 34 '0.04713'
$ws.Cells.Item(34, 5).Value = '  -3.48%  '

# Row 35: ImmutableX
$ws.Cells.Item(35, 2).Value = 'ImmutableX'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-PriceText 35 '0.7230'
$ws.Cells.Item(35, 5).Value = '  -3.07%  '

# Row 36: ARBITRUM
$ws.Cells.Item(36, 2).Value = 'ARBITRUM'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-PriceText 36 '1.099'
$ws.Cells.Item(36, 5).Value = '  -5.61%  '

# Row 37: Frax
$ws.Cells.Item(37, 2).Value = 'Frax'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-PriceText 37 '1.000'
$ws.Cells.Item(37, 5).Value = '  -0.04%  '

# Row 38: HuobiToken
$ws.Cells.Item(38, 2).Value = 'HuobiToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-PriceText 38 '2.717'
$ws.Cells.Item(38, 5).Value = '  -0.69%  '

# Row 39: VeChain
$ws.Cells.Item(39, 2).Value = 'VeChain'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-PriceText 39 '0.01897'
$ws.Cells.Item(39, 5).Value = '  -5.02%  '

# Row 40: MXToken
$ws.Cells.Item(40, 2).Value = 'MXToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-PriceText 40 '2.591'
$ws.Cells.Item(40, 5).Value = '  -3.63%  '

# Row 41: FraxShare
$ws.Cells.Item(41, 2).Value = 'FraxShare'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-PriceText 41 '6.212'
$ws.Cells.Item(41, 5).Value = '  -4.62%  '

# Row 42: Aave
$ws.Cells.Item(42, 2).Value = 'Aave'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-PriceText 42 '72.38'
$ws.Cells.Item(42, 5).Value = '  -6.77%  '

# Row 43: RenderToken
$ws.Cells.Item(43, 2).Value = 'RenderToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-PriceText 43 '1.971'
$ws.Cells.Item(43, 5).Value = '  -7.03%  '

# Row 44: Quant
$ws.Cells.Item(44, 2).Value = 'Quant'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-PriceText 44 '106.20'
$ws.Cells.Item(44, 5).Value = '  -2.58%  '

# Row 45: TrustWalletToken
$ws.Cells.Item(45, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-PriceText 45 '0.8555'
$ws.Cells.Item(45, 5).Value = '  -5.16%  '

# Row 46: Maker
Set-PriceText 46 '1.043.51'
$ws.Cells.Item(46, 5).Value = '  +5.08%  '

# Row 47: PaxDollar
$ws.Cells.Item(47, 2).Value = 'PaxDollar'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-PriceText 47 '1.002'
$ws.Cells.Item(47, 5).Value = '  +0.02%  '

# Row 48: TheSandbox
$ws.Cells.Item(48, 2).Value = 'TheSandbox'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-PriceText 48 '0.4201'
$ws.Cells.Item(48, 5).Value = '  -5.35%  '

# Row 49: Aptos
$ws.Cells.Item(49, 2).Value = 'Aptos'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-PriceText 49 '7.377'
$ws.Cells.Item(49, 5).Value = '  -7.85%  '

# Row 50: Algorand
$ws.Cells.Item(50, 2).Value = 'Algorand'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-PriceText 50 '0.1186'
$ws.Cells.Item(50, 5).Value = '  -4.59%  '

# Row 51: Elrond
$ws.Cells.Item(51, 2).Value = 'Elrond'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-PriceText 51 '34.38'
$ws.Cells.Item(51, 5).Value = '  -4.16%  '

$helper.Clear()
$excel.CutCopyMode = $false
